$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking strings
# (e.g. "1.00", "206.91") are not silently converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '67.936.53'
$ws.Range("E2").Value = '  +2.18%  '

$ws.Range("D3").Value = '3.564.70'
$ws.Range("E3").Value = '  +0.29%  '

$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.17%  '

$ws.Range("D5").Value = '206.91'
$ws.Range("E5").Value = '  +11.39%  '

$ws.Range("D6").Value = '562.40'
$ws.Range("E6").Value = '  -1.55%  '

$ws.Range("D7").Value = '3.562.41'
$ws.Range("E7").Value = '  +0.35%  '

$ws.Range("D8").Value = '0.609'
$ws.Range("E8").Value = '  -0.71%  '

$ws.Range("E9").Value = '  -0.07%  '

$ws.Range("D10").Value = '0.675'
$ws.Range("E10").Value = '  +0.73%  '

$ws.Range("D11").Value = '60.51'
$ws.Range("E11").Value = '  +9.84%  '

$ws.Range("D12").Value = '0.146'
$ws.Range("E12").Value = '  -2.43%  '

$ws.Range("D13").Value = '0.0000278'
$ws.Range("E13").Value = '  +6.33%  '

$ws.Range("D14").Value = '10.17'
$ws.Range("E14").Value = '  +3.82%  '

$ws.Range("D15").Value = '4.149.70'
$ws.Range("E15").Value = '  +0.71%  '

$ws.Range("D16").Value = '3.576.31'
$ws.Range("E16").Value = '  +0.70%  '

$ws.Range("E17").Value = '  +0.72%  '

$ws.Range("D18").Value = '18.75'
$ws.Range("E18").Value = '  +2.85%  '

$ws.Range("D19").Value = '67.803.24'
$ws.Range("E19").Value = '  +2.09%  '

$ws.Range("D20").Value = '12.13'
$ws.Range("E20").Value = '  +0.74%  '

$ws.Range("D21").Value = '1.05'
$ws.Range("E21").Value = '  -0.20%  '

$ws.Range("D22").Value = '399.96'
$ws.Range("E22").Value = '  +3.11%  '

$ws.Range("D23").Value = '12.43'
$ws.Range("E23").Value = '  +11.13%  '

$ws.Range("D24").Value = '4.13'
$ws.Range("E24").Value = '  -0.96%  '

$ws.Range("D25").Value = '83.94'
$ws.Range("E25").Value = '  -1.50%  '

$ws.Range("D26").Value = '2.85'
$ws.Range("E26").Value = '  -1.75%  '

$ws.Range("D27").Value = '12.34'
$ws.Range("E27").Value = '  -0.13%  '

$ws.Range("D28").Value = '3.86'
$ws.Range("E28").Value = '  +8.46%  '

$ws.Range("D29").Value = '9.14'
$ws.Range("E29").Value = '  +3.61%  '

$ws.Range("D30").Value = '7.67'
$ws.Range("E30").Value = '  +1.67%  '

$ws.Range("D31").Value = '31.28'
$ws.Range("E31").Value = '  +1.64%  '

$ws.Range("D32").Value = '660.19'
$ws.Range("E32").Value = '  +4.48%  '

$ws.Range("D33").Value = '12.00'
$ws.Range("E33").Value = '  -0.93%  '

$ws.Range("D34").Value = '63.07'
$ws.Range("E34").Value = '  -0.26%  '

$ws.Range("E35").Value = '  -1.04%  '

$ws.Range("D36").Value = '40.85'
$ws.Range("E36").Value = '  -2.07%  '

$ws.Range("D37").Value = '0.406'
$ws.Range("E37").Value = '  +1.13%  '

$ws.Range("E38").Value = '  -0.15%  '

$ws.Range("E39").Value = '  +11.94%  '

$ws.Range("D40").Value = '0.0₃0745'
$ws.Range("E40").Value = '  +0.02%  '

$ws.Range("D41").Value = '3.154.61'
$ws.Range("E41").Value = '  +0.70%  '

$ws.Range("E42").Value = '  +0.06%  '

$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  +0.11%  '

$ws.Range("D44").Value = '2.67'
$ws.Range("E44").Value = '  +0.83%  '

$ws.Range("D45").Value = '2.78'
$ws.Range("E45").Value = '  +11.47%  '

$ws.Range("D46").Value = '0.0407'
$ws.Range("E46").Value = '  -0.69%  '

$ws.Range("D47").Value = '0.129'
$ws.Range("E47").Value = '  -0.15%  '

$ws.Range("B48").Value = 'ApeXProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D48").Value = '3.06'
$ws.Range("E48").Value = '  +0.01%  '

$ws.Range("B49").Value = 'THORChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D49").Value = '8.61'
$ws.Range("E49").Value = '  +2.67%  '

$ws.Range("B50").Value = 'dogwifhat'
$ws.Range("C50").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D50").Value = '2.56'
$ws.Range("E50").Value = '  +8.88%  '

$ws.Range("D51").Value = '138.29'
$ws.Range("E51").Value = '  -0.76%  '

# Restore the original (default/general) formatting on column D so the
# cell style matches the source workbook (no explicit style index).
$ws.Range("D2:D51").ClearFormats()
